$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.118.55"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "2.524.92"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'537.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").Value = "'137.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.79%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("D9").Value = "2.523.76"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").Value = "'5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("E13").Value = "  -2.97%  "

$ws.Range("D14").Value = "2.956.45"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").Value = "'23.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("D16").Value = "59.006.57"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").Value = "2.533.33"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").Value = "'323.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  +1.72%  "

$ws.Range("D24").Value = "'65.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.15%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -2.11%  "

$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").Value = "'7.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.74%  "

$ws.Range("D29").Value = "0.0₃0773"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").Value = "'6.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.29%  "

$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("D32").Value = "'167.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.84%  "

$ws.Range("E33").Value = "  +5.40%  "

$ws.Range("E35").Value = "  +1.83%  "

$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("D37").Value = "'4.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.88%  "

$ws.Range("D38").Value = "'1.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.44%  "

$ws.Range("D39").Value = "'36.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").Value = "'3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.83%  "

$ws.Range("D42").Value = "'284.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("D43").Value = "'5.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("D44").Value = "'132.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.86%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("E46").Value = "  +1.53%  "

$ws.Range("D47").Value = "'10.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'0.0925"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.31%  "

$ws.Range("D49").Value = "'0.0508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("E50").Value = "  -1.82%  "

$ws.Range("D51").Value = "'17.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.96%  "
